# Commit: "added github link to pptx"
# Adds a centered "GITHUB [ https://github.com/roberto-chan/bicep-intro ]"
# textbox (with the URL hyperlinked) to the bottom of slide 2 ("Agenda" /
# demo slide), just below the two pictures already on that slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Slide is 12192000 x 6858000 EMU (960 x 540 pt); Shapes.AddTextbox takes
# points, so convert the target EMU box (0,6256800)-(12115800,457200) by /12700.
$left   = 0 / 12700.0
$top    = 6256800 / 12700.0
$width  = 12115800 / 12700.0
$height = 457200 / 12700.0

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = ""

$shp.Fill.Visible = $false
$shp.Line.Visible = $false
$shp.Line.Weight = 0

$tf = $shp.TextFrame
$tf.MarginLeft = 90000 / 12700.0
$tf.MarginRight = 90000 / 12700.0
$tf.MarginTop = 45000 / 12700.0
$tf.MarginBottom = 45000 / 12700.0
$tf.AutoSize = 0
$tf.VerticalAnchor = 1

$part1 = "GITHUB"
$part2 = " [ "
$part3 = "https://github.com/roberto-chan/bicep-intro"
$part4 = " ]"

$tr = $tf.TextRange
$tr.Text = $part1 + $part2 + $part3 + $part4
$tr.ParagraphFormat.Alignment = 2

# Shared baseline formatting for the whole run of text.
$tr.Font.Bold = $true
$tr.Font.Size = 22
$tr.Font.Name = "Arial"
$tr.Font.Color.RGB = 0
$tr.Font.Spacing = -0.01
$tr.Font.Strikethrough = $false

$start1 = 1
$start2 = $start1 + $part1.Length
$start3 = $start2 + $part2.Length
$start4 = $start3 + $part3.Length

$run1 = $tr.Characters($start1, $part1.Length)
$run1.Font.Bold = $true
$run1.Font.Size = 22
$run1.Font.Name = "Arial"
$run1.Font.Spacing = -0.01
$run1.Font.Strikethrough = $false
$run1.Font.Color.RGB = 6710886

$run2 = $tr.Characters($start2, $part2.Length)
$run2.Font.Bold = $true
$run2.Font.Size = 22
$run2.Font.Name = "Arial"
$run2.Font.Spacing = -0.01
$run2.Font.Strikethrough = $false
$run2.Font.Color.RGB = 0

$run3 = $tr.Characters($start3, $part3.Length)
$run3.Font.Bold = $true
$run3.Font.Size = 22
$run3.Font.Name = "Arial"
$run3.Font.Spacing = -0.01
$run3.Font.Strikethrough = $false
$run3.Font.Color.RGB = 0
$run3.ActionSettings.Item(1).Hyperlink.Address = $part3

$run4 = $tr.Characters($start4, $part4.Length)
$run4.Font.Bold = $true
$run4.Font.Size = 22
$run4.Font.Name = "Arial"
$run4.Font.Spacing = -0.01
$run4.Font.Strikethrough = $false
$run4.Font.Color.RGB = 0

Write-Output "added github link shape to slide 2"
